# Documentation.docx - add basic descriptions for every GUI mock-up element.
$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$vNs = 'xmlns:v="urn:schemas-microsoft-com:vml"'
$oNs = 'xmlns:o="urn:schemas-microsoft-com:office:office"'
$rNs = 'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

# ---------------------------------------------------------------------------
# 1) Remove the trailing "Networking" Heading3 paragraph (it gets moved into
#    the body text earlier in the document, right after the mock-up notes).
# ---------------------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
if ($last.Range.Text.TrimEnd([char]13, [char]7) -eq "Networking") {
    $last.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) "Technical I" + _GoBack bookmark + "nformation" -> single run, bookmark
#    removed here (it gets re-created further up, in the new Raw/Processed
#    Mode paragraph).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Technical Information") {
        $xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr><w:r><w:t>Technical Information</w:t></w:r></w:p>"
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the "(explanation goes here)" placeholder paragraph with the
#    real descriptions of every mock-up element, plus the "Networking"
#    heading lead-in paragraph (with the relocated _GoBack bookmark and a
#    lastRenderedPageBreak marker on "Networking").
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "(explanation goes here)") {
        $xml = @"
<w:p $wNs><w:r><w:t>Object Display Canvas:</w:t></w:r><w:r><w:t xml:space="preserve"> Either displays the greyscale image in Raw mode or displays the detected objects in Processed mode. The greyscale image and detected objects will be scaled according to the current camera settings. The detected objects will be color coded according to their possible danger to the SRS. In Processed mode, objects can be selected and their details can be viewed from the Selected Object Information Panel. </w:t></w:r></w:p><w:p $wNs><w:r><w:t>Selected Object Information Panel:</w:t></w:r><w:r><w:t xml:space="preserve"> Displays the</w:t></w:r><w:r><w:t xml:space="preserve"> current selected</w:t></w:r><w:r><w:t xml:space="preserve"> object's camera assigned ID, estimated size, and estimated velocity</w:t></w:r><w:r><w:t xml:space="preserve"> while in Processed Mode</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p $wNs><w:r><w:t>Camera Zoom:</w:t></w:r><w:r><w:t xml:space="preserve"> Adjusts the camera's zoom toward the center of its field of view.</w:t></w:r></w:p><w:p $wNs><w:r><w:t>Section Size:</w:t></w:r><w:r><w:t xml:space="preserve"> Adjusts the image processing section size.</w:t></w:r></w:p><w:p $wNs><w:r><w:t>Section Overlap:</w:t></w:r><w:r><w:t xml:space="preserve"> Adjusts the amount of overlap each section is allowed to have with its neighbors.</w:t></w:r></w:p><w:p $wNs><w:r><w:t>Raw/Processed Mode:</w:t></w:r><w:r><w:t xml:space="preserve"> Raw Mode returns the greyscale image with no processing. Processed Mode returns all objects detected processed and its relative coordinates on the field of view.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p $wNs><w:r><w:lastRenderedPageBreak/><w:t>Networking</w:t></w:r></w:p>
"@
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Caption paragraph: drop the spell-check markers around "SpaceRock" and
#    merge the trailing runs into a single run, keeping the SEQ field intact.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Figure ") -and $p.Range.Text.Contains("Mock-Up of SpaceRock GUI")) {
        $xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Caption`"/><w:jc w:val=`"center`"/></w:pPr><w:r><w:t xml:space=`"preserve`">Figure </w:t></w:r><w:r><w:fldChar w:fldCharType=`"begin`"/></w:r><w:r><w:instrText xml:space=`"preserve`"> SEQ Figure \* ARABIC </w:instrText></w:r><w:r><w:fldChar w:fldCharType=`"separate`"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r><w:r><w:fldChar w:fldCharType=`"end`"/></w:r><w:r><w:t>: Mock-Up of SpaceRock GUI</w:t></w:r></w:p>"
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 5) Embedded Visio OLE object: refreshed shape id / size / ObjectID.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Trim([char]13, [char]7) -eq "" -and $p.Range.Fields.Count -ge 1) {
        $xml = @"
<w:p $wNs $vNs $oNs $rNs><w:pPr><w:keepNext/><w:jc w:val="center"/></w:pPr><w:r><w:object w:dxaOrig="8846" w:dyaOrig="5920"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:442.5pt;height:296.25pt" o:ole=""><v:imagedata r:id="rId4" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Visio.Drawing.15" ShapeID="_x0000_i1025" DrawAspect="Content" ObjectID="_1547901003" r:id="rId5"/></w:object></w:r></w:p>
"@
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 6) Title: drop the spell-check markers around "SpaceRock" and merge into a
#    single run.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading1`"/></w:pPr><w:r><w:t>SpaceRock GUI Documentation</w:t></w:r></w:p>"
$p1.Range.InsertXML($xml)
